# UC007 - Listar Autorizações de Pagamento Pendentes
# v1.1.1 -> v1.2
#
# The second step of TC3 and the second step of TC4 were swapped:
#   TC3 step 2 (row 28) used to describe the "filter by user" action/result;
#   it now gets TC4's "realizar a autorização de pagamento" action/result.
#   TC4 step 2 (row 36) used to describe the "realizar a autorização de pagamento" action/result;
#   it now gets TC3's "filter by user" action/result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC3 block - row 28 (2nd step): was "filter by user", now "realizar a autorização de pagamento"
$ws.Range("B28").Value = "Chefe Clica para realizar a autorização de pagamento."
$ws.Range("D28").Value = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"

# TC4 block - row 36 (2nd step): was "realizar a autorização de pagamento", now "filter by user"
$ws.Range("B36").Value = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$ws.Range("D36").Value = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."
